$d = $word.ActiveDocument

# 1. Rename the first referenced document from the old MOA/AOA pdf to the
#    resume pdf (heading2 paragraph just after the report title).
$d.Content.Find.Execute("GEPL-MOA-and-AOA-August.pdf", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Thatikonda Narendra_AI_Developer_3Yrs+_Exp (1).pdf", 2)

# 2. Collapse the "Document Type: MOA" / "Overall Risk: Medium" /
#    "High: 0, Medium: 4, Low: 1, Total: 5" / second-document heading /
#    "Document Type: UNKNOWN" block down into a single
#    "Document Type: LOAN_AGREEMENT" paragraph, leaving the trailing
#    "Overall Risk: No Risk" / "High: 0, Medium: 0, Low: 0, Total: 0"
#    paragraphs untouched.
$pDocType = $d.Paragraphs.Item(3)
$pOverallRiskMedium = $d.Paragraphs.Item(4)
$pCounts = $d.Paragraphs.Item(5)
$pSecondHeading = $d.Paragraphs.Item(6)
$pUnknown = $d.Paragraphs.Item(7)

$pDocType.Range.Text = "Document Type: LOAN_AGREEMENT"

$toRemove = $d.Range($pOverallRiskMedium.Range.Start, $pUnknown.Range.End)
$toRemove.Delete()
